# Add character styles GaNStyle, GaNParagraph, GaNLinks and apply them
# to the relevant runs, per the "Add styles to the new paragraphs" commit.

$d = $word.ActiveDocument

# --- Create the three new character styles -------------------------------

$GaNStyle = $d.Styles.Add("GaNStyle", 2)
$GaNStyle.Font.Name = "Calibri"
$GaNStyle.Font.Size = 14

$GaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$GaNParagraph.Font.Name = "Calibri"
$GaNParagraph.Font.Size = 10

$GaNLinks = $d.Styles.Add("GaNLinks", 2)
$GaNLinks.Font.Name = "Calibri"
$GaNLinks.Font.Size = 9.5
$GaNLinks.Font.Bold = $true
$GaNLinks.Font.Color = 8388608
$GaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas da campaña..." run -------------------

$datesText = "Datas da campaña de 2022 que usan Constelación de Leo: 14-23 de abril, 14-23 de maio"
$range = $d.Content
$found = $range.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $found = $range.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Grazas por participar..." run ------------

$thanksText = "Grazas por participar nesta campaña global de medida da contaminación lumínica mediante a observación das estrelas máis febles que podes albiscar. Localizando e observando a  Constelación de Leo e comparándoa co que aparece nos mapas estelares recollidos neste documento podes saber canto contribúen á contaminación lumínica os sistemas de iluminación que hai no teu barrio ou vila. As túas achegas á base de datos en liña de GLOBE at Night (O MUNDO á Noite) servirán para documentar a calidade do ceo nocturno."
$range2 = $d.Content
$found2 = $range2.Find.Execute($thanksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found2) {
    $range2.Style = "GaNParagraph"
    $range2.Collapse(0)
    $found2 = $range2.Find.Execute($thanksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNLinks to the "Os mapas de estrelas..." run -----------------

$mapsText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$found3 = $range3.Find.Execute($mapsText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found3) {
    $range3.Style = "GaNLinks"
    $range3.Collapse(0)
    $found3 = $range3.Find.Execute($mapsText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Output "Styles added and applied."
